# Daily attendance processing - 2025-12-30 18:09:29
# Normalize the "Recorded By" (column G) ordering: swap
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# for every row where it appears, leaving any other values
# (e.g. lone "System" or lone "dnasr281@gmail.com") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Recorded By cells updated: $updated"
